$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ASGS SA1/SA2 download rows (appended below the existing MB rows)
$url_sa1 = "https://www.abs.gov.au/ausstats/subscriber.nsf/log?openagent&1270055004_sa1_ucl_sosr_sos_2016_aust_csv.zip&1270.0.55.004&Data%20Cubes&EE5F4698A91AD2F8CA2581B1000E09B0&0&July%202016&09.10.2017&Latest"
$url_sa2 = "https://www.abs.gov.au/ausstats/subscriber.nsf/log?openagent&1270055004_sa2_sua_2016_aust_csv.zip&1270.0.55.004&Data%20Cubes&D6E51168BD6DC248CA2581B1000E0A48&0&July%202016&09.10.2017&Latest"
$name_sa1 = "1270055004_sa1_ucl_sosr_sos_2016_aust_csv.zip"
$name_sa2 = "1270055004_sa2_sua_2016_aust_csv.zip"

# Column A (urls) first, top to bottom
$ws.Range("A11").Value = $url_sa1
$ws.Range("A12").Value = $url_sa2

# Column B (file names), filled in B12 then B11 order
$ws.Range("B12").Value = $name_sa2
$ws.Range("B11").Value = $name_sa1

# Hyperlinks: A12 added before A11
$ws.Hyperlinks.Add($ws.Range("A12"), $url_sa2) | Out-Null
$ws.Range("A12").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("A11"), $url_sa1) | Out-Null
$ws.Range("A11").Style = "Hyperlink"

# Widen columns to fit the new, longer content
$ws.Columns("A").ColumnWidth = 200.3
$ws.Columns("B").ColumnWidth = 44.0

$ws.Range("B19").Select() | Out-Null
